# Update "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN..is_holiday_week right by one)
#  - shorten the week labels (W01 -> W1, ... W16 -> W16)
#  - populate the new Week_Start_Date column with the week's start date (as text)
#  - correct a handful of MyForecast values
#  - store is_holiday_week as a boolean instead of a number
# Then update the "Summary" sheet's forecast totals to match the corrected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column B ("ASIN" and everything right of it shifts over to make room
# for the new "Week_Start_Date" column).
$ws.Columns("B:B").Insert()

# Header row
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Row data: Week, Week_Start_Date, MyForecast, AmazonMean, AmazonP70, AmazonP80, AmazonP90
$rows = @(
    @{ r = 2;  week = "W1";  date = "2025-01-05"; myf = 21; mean = 19; p70 = 22; p80 = 26; p90 = 31 },
    @{ r = 3;  week = "W2";  date = "2025-01-12"; myf = 21; mean = 20; p70 = 24; p80 = 28; p90 = 35 },
    @{ r = 4;  week = "W3";  date = "2025-01-19"; myf = 23; mean = 21; p70 = 25; p80 = 30; p90 = 37 },
    @{ r = 5;  week = "W4";  date = "2025-01-26"; myf = 22; mean = 21; p70 = 26; p80 = 31; p90 = 39 },
    @{ r = 6;  week = "W5";  date = "2025-02-02"; myf = 24; mean = 22; p70 = 27; p80 = 32; p90 = 40 },
    @{ r = 7;  week = "W6";  date = "2025-02-09"; myf = 23; mean = 21; p70 = 26; p80 = 31; p90 = 39 },
    @{ r = 8;  week = "W7";  date = "2025-02-16"; myf = 25; mean = 22; p70 = 27; p80 = 33; p90 = 43 },
    @{ r = 9;  week = "W8";  date = "2025-02-23"; myf = 24; mean = 23; p70 = 28; p80 = 34; p90 = 44 },
    @{ r = 10; week = "W9";  date = "2025-03-02"; myf = 25; mean = 22; p70 = 27; p80 = 34; p90 = 44 },
    @{ r = 11; week = "W10"; date = "2025-03-09"; myf = 24; mean = 22; p70 = 27; p80 = 34; p90 = 45 },
    @{ r = 12; week = "W11"; date = "2025-03-16"; myf = 26; mean = 23; p70 = 28; p80 = 35; p90 = 47 },
    @{ r = 13; week = "W12"; date = "2025-03-23"; myf = 25; mean = 23; p70 = 29; p80 = 36; p90 = 48 },
    @{ r = 14; week = "W13"; date = "2025-03-30"; myf = 26; mean = 22; p70 = 27; p80 = 34; p90 = 45 },
    @{ r = 15; week = "W14"; date = "2025-04-06"; myf = 25; mean = 21; p70 = 26; p80 = 34; p90 = 45 },
    @{ r = 16; week = "W15"; date = "2025-04-13"; myf = 27; mean = 22; p70 = 27; p80 = 34; p90 = 47 },
    @{ r = 17; week = "W16"; date = "2025-04-20"; myf = 26; mean = 22; p70 = 27; p80 = 34; p90 = 46 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.week
    # Leading apostrophe forces the date-look-alike string to stay text instead of
    # being auto-converted into a date serial number.
    $ws.Cells.Item($r, 2).Value = "'" + $row.date
    $ws.Cells.Item($r, 4).Value = $row.myf
    $ws.Cells.Item($r, 5).Value = $row.mean
    $ws.Cells.Item($r, 6).Value = $row.p70
    $ws.Cells.Item($r, 7).Value = $row.p80
    $ws.Cells.Item($r, 8).Value = $row.p90
    # is_holiday_week becomes a real boolean (FALSE) instead of numeric 0
    $ws.Cells.Item($r, 10).Value = $false
}

# Update the Summary sheet's forecast totals to reflect the corrected MyForecast numbers.
# (kept as text, like the rest of the Value column on this sheet - leading apostrophe
# keeps the numeric-looking string from being auto-converted into a number)
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(9, 2).Value = "'389"
$summary.Cells.Item(10, 2).Value = "'184"
$summary.Cells.Item(11, 2).Value = "'88"
